$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110..182 down to 111..183
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new data record
$ws.Cells.Item(110, 1).Value = 3
$ws.Cells.Item(110, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(110, 3).Value = "Coquimbo"
$ws.Cells.Item(110, 4).Value = 44582
$ws.Cells.Item(110, 5).Value = 5
$ws.Cells.Item(110, 6).Value = "Fruta"
$ws.Cells.Item(110, 7).Value = 100101
$ws.Cells.Item(110, 8).Value = "Berries"
$ws.Cells.Item(110, 9).Value = 100101001
$ws.Cells.Item(110, 10).Value = "Arándano (blue)"
$ws.Cells.Item(110, 11).Value = "Sin especificar"
$ws.Cells.Item(110, 12).Value = "Primera"
$ws.Cells.Item(110, 13).Value = 160
$ws.Cells.Item(110, 14).Value = 4000
$ws.Cells.Item(110, 15).Value = 4500
$ws.Cells.Item(110, 16).Value = 4250
$ws.Cells.Item(110, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(110, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(110, 19).Value = 2125
$ws.Cells.Item(110, 20).Value = 2
